$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 94 currently carries the "last row" date style (s=3, YYYY-MM-DD).
# Once we append a new row it becomes a regular data row again, so restore
# the normal date style (s=2, YYYY-MM-DD HH:MM:SS) used by the other rows.
$ws.Range("A94").NumberFormat = $ws.Range("A93").NumberFormat

# Append the new day's data as row 95, giving the date cell the "last row"
# style that A94 used to have.
$ws.Range("A95").Value = 45834
$ws.Range("A95").NumberFormat = "YYYY-MM-DD"

$ws.Range("B95").Value = 403
$ws.Range("C95").Value = 399
$ws.Range("D95").Value = 411
